$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.821.58"
$ws.Range("E2").Value2 = "  -1.59%  "
$ws.Range("D3").Value2 = "1.889.24"
$ws.Range("E3").Value2 = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.7726"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -5.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "244.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3126"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "25.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -7.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07225"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -0.67%  "
$ws.Range("E11").Value2 = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.7657"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -3.70%  "
$ws.Range("B13").Value2 = "Polkadot"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.507"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +1.76%  "
$ws.Range("B14").Value2 = "WrappedEther"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value2 = "1.908.08"
$ws.Range("E14").Value2 = "  -0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "92.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "6.167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +1.44%  "
$ws.Range("D17").Value2 = "29.825.58"
$ws.Range("E17").Value2 = "  -1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "13.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "243.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -3.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "0.000007776"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -1.09%  "
$ws.Range("B21").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value2 = "2.162.76"
$ws.Range("E21").Value2 = "  -1.14%  "
$ws.Range("B22").Value2 = "Dai"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "8.179"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.1573"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -6.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "9.431"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "162.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -3.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "18.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.037"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -5.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.448"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +5.57%  "
$ws.Range("E31").Value2 = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.450"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.080"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05504"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.257"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -3.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.7488"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +0.27%  "
$ws.Range("E37").Value2 = "  +0.31%  "
$ws.Range("E38").Value2 = "  -3.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value2 = "  -1.32%  "
$ws.Range("D41").Value2 = "1.145.25"
$ws.Range("E41").Value2 = "  +10.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "73.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.4420"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "5.906"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.8490"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -0.82%  "
$ws.Range("E46").Value2 = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "102.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.883"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "9.904"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "7.441"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "3.007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -3.69%  "
